$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($RangeAddr, $Val)
    $r = $ws.Range($RangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '37.000.15'
Set-TextValue 'E2' '  +1.19%  '

# Row 3
Set-TextValue 'D3' '1.980.85'
Set-TextValue 'E3' '  +1.03%  '

# Row 5
Set-TextValue 'D5' '245.59'
Set-TextValue 'E5' '  +0.46%  '

# Row 6
Set-TextValue 'E6' '  +1.83%  '

# Row 7
Set-TextValue 'D7' '61.25'
Set-TextValue 'E7' '  +4.11%  '

# Row 8
Set-TextValue 'E8' '  +0.03%  '

# Row 9
Set-TextValue 'E9' '  +1.76%  '

# Row 10
Set-TextValue 'D10' '0.0799'
Set-TextValue 'E10' '  -1.44%  '

# Row 11
Set-TextValue 'E11' '  +0.00%  '

# Row 12
Set-TextValue 'D12' '15.00'
Set-TextValue 'E12' '  +9.36%  '

# Row 13
Set-TextValue 'B13' 'Avalanche'
Set-TextValue 'C13' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D13' '22.19'
Set-TextValue 'E13' '  -0.28%  '

# Row 14
Set-TextValue 'B14' 'Polygon'
Set-TextValue 'C14' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.844'
Set-TextValue 'E14' '  +1.99%  '

# Row 15
Set-TextValue 'D15' '2.273.45'
Set-TextValue 'E15' '  +1.10%  '

# Row 16
Set-TextValue 'D16' '5.48'
Set-TextValue 'E16' '  +3.82%  '

# Row 17
Set-TextValue 'D17' '1.983.29'
Set-TextValue 'E17' '  +1.27%  '

# Row 18
Set-TextValue 'D18' '36.896.45'
Set-TextValue 'E18' '  +1.08%  '

# Row 19
Set-TextValue 'D19' '70.13'
Set-TextValue 'E19' '  +0.35%  '

# Row 20
Set-TextValue 'E20' '  +0.24%  '

# Row 21
Set-TextValue 'E21' '  +1.90%  '

# Row 22
Set-TextValue 'D22' '230.23'
Set-TextValue 'E22' '  +0.61%  '

# Row 23
Set-TextValue 'E23' '  +0.01%  '

# Row 24
Set-TextValue 'D24' '2.50'
Set-TextValue 'E24' '  +1.19%  '

# Row 25
Set-TextValue 'E25' '  +0.13%  '

# Row 26
Set-TextValue 'E26' '  +8.81%  '

# Row 27
Set-TextValue 'D27' '9.28'
Set-TextValue 'E27' '  +0.68%  '

# Row 28
Set-TextValue 'D28' '163.28'
Set-TextValue 'E28' '  +1.97%  '

# Row 29
Set-TextValue 'D29' '19.56'
Set-TextValue 'E29' '  +0.54%  '

# Row 30
Set-TextValue 'E30' '  +17.67%  '

# Row 31
Set-TextValue 'E31' '  +1.89%  '

# Row 32
Set-TextValue 'D32' '4.86'
Set-TextValue 'E32' '  +2.77%  '

# Row 33
Set-TextValue 'E33' '  +0.20%  '

# Row 34
Set-TextValue 'D34' '4.52'
Set-TextValue 'E34' '  +5.23%  '

# Row 35
Set-TextValue 'D35' '2.32'
Set-TextValue 'E35' '  +2.68%  '

# Row 36
Set-TextValue 'E36' '  -0.09%  '

# Row 37
Set-TextValue 'D37' '3.36'
Set-TextValue 'E37' '  -0.15%  '

# Row 38
Set-TextValue 'E38' '  +0.48%  '

# Row 39
Set-TextValue 'D39' '5.52'
Set-TextValue 'E39' '  -7.38%  '

# Row 40
Set-TextValue 'D40' '0.0976'
Set-TextValue 'E40' '  -0.99%  '

# Row 41
Set-TextValue 'E41' '  +1.01%  '

# Row 42
Set-TextValue 'E42' '  +0.27%  '

# Row 43
Set-TextValue 'E43' '  +0.79%  '

# Row 44
Set-TextValue 'E44' '  +2.82%  '

# Row 45
Set-TextValue 'B45' 'Maker'
Set-TextValue 'C45' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D45' '1.372.31'
Set-TextValue 'E45' '  +0.49%  '

# Row 46
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '90.16'
Set-TextValue 'E46' '  +2.68%  '

# Row 47
Set-TextValue 'E47' '  +0.17%  '

# Row 48
Set-TextValue 'E48' '  +1.35%  '

# Row 49
Set-TextValue 'E49' '  -0.49%  '

# Row 50
Set-TextValue 'E50' '  +6.00%  '

# Row 51
Set-TextValue 'D51' '1.95'
Set-TextValue 'E51' '  +9.12%  '
